# Automatski commit: 2025-02-09 22:55
# Update report number 07/08 -> 16, month november -> februar, year 2024 -> 2025,
# and add a page break + new (left aligned) heading paragraph after the date line.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Body: "08_lav03_vaja" heading -> "16_lav03_vaja"
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs(3)
$p3.Range.Find.Execute("08", $true, $false, $false, $false, $false, $true, 1, $false, "16", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Body: "Vaja 8" -> "Vaja 16"
# ---------------------------------------------------------------------------
$p5 = $d.Paragraphs(5)
$p5.Range.Find.Execute("8", $true, $false, $false, $false, $false, $true, 1, $false, "16", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Body: "Ljubljana, november 2024" -> "Ljubljana, februar 2025"
# ---------------------------------------------------------------------------
$p7 = $d.Paragraphs(7)
$p7.Range.Find.Execute("november", $true, $false, $false, $false, $false, $true, 1, $false, "februar", 2) | Out-Null
$p7 = $d.Paragraphs(7)
$p7.Range.Find.Execute("2024", $true, $false, $false, $false, $false, $true, 1, $false, "2025", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) Body: insert a page break paragraph followed by a new, empty
#    "Naslov21" (left aligned) paragraph right after the date line.
# ---------------------------------------------------------------------------
$p7 = $d.Paragraphs(7)
$rEnd = $d.Range($p7.Range.End, $p7.Range.End)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
'<pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part>' +
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' +
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
'<w:body>' +
'<w:p><w:r><w:br w:type="page"/></w:r></w:p>' +
'<w:p><w:pPr><w:pStyle w:val="Naslov21"/><w:jc w:val="left"/></w:pPr></w:p>' +
'</w:body></w:document>' +
'</pkg:xmlData></pkg:part></pkg:package>'

$rEnd.InsertXML($xml) | Out-Null

# ---------------------------------------------------------------------------
# 5) header1.xml (default header): "poročilo 07_lav03_vaja" -> "poročilo 16_lav03_vaja"
# ---------------------------------------------------------------------------
$sec = $d.Sections(1)
$hdr1 = $sec.Headers(1)
$hdr1.Range.Find.Execute("poročilo 07", $true, $false, $false, $false, $false, $true, 1, $false, "poročilo 16", 2) | Out-Null

# ---------------------------------------------------------------------------
# 6) header2.xml (first-page header): "Poročilo 08_lav03_vaja" -> "Poročilo 16_lav03_vaja"
# ---------------------------------------------------------------------------
$hdr2 = $sec.Headers(2)
$hdr2.Range.Find.Execute("Poročilo 08", $true, $false, $false, $false, $false, $true, 1, $false, "Poročilo 16", 2) | Out-Null

Write-Host "Done"
